$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "38.287.99"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +3.75%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.064.01"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("E4").Value = "  +0.48%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "230.93"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.09%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.617"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.30%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "58.74"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +8.24%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +3.78%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0811"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +4.42%  "
$ws.Range("E11").Value = "  +0.83%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "2.367.62"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +3.30%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "14.71"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +5.08%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "20.75"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.757"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("E16").Value = "  +4.76%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.055.87"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.30%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "38.059.94"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +3.38%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.18"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.98%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "69.98"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("E21").Value = "  +2.98%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "225.21"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  +5.10%  "
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("E28").Value = "  +8.87%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "19.11"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("E30").Value = "  +2.45%  "
$ws.Range("E31").Value = "  +2.82%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.57"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +2.32%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.64"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +6.18%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0616"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.74%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +7.71%  "
$ws.Range("E36").Value = "  +1.11%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "6.06"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +16.07%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "3.34"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +7.02%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  +3.13%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "98.67"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.484.27"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.56%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0951"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.58%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "16.95"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +5.16%  "
$ws.Range("E45").Value = "  +3.94%  "
$ws.Range("E46").Value = "  +1.35%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.07"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +18.18%  "
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("E50").Value = "  +0.02%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.252.22"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +3.32%  "
